$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C (rows 2 through 110) holds a "last updated" date serial that
# advances by one day (45179 -> 45180) for every data row in the sheet.
$ws.Range("C2:C110").Value = 45180
